# Deb FG Update Template How-To - "save docs to pdf"
#
# The template had a stray blank row at the very top (row 1) and a couple of
# stray blank rows further down (rows 16-17, between the "senario" note row
# and the second table). Both of these were removed, which shifts all of the
# real content up and makes the sheet tidy/printable; the two tables
# automatically re-anchor to the new (smaller) ranges. The sheet was also set
# to print in landscape orientation (prep for exporting/printing to PDF).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank spacer row at the very top of the sheet - everything
# below (including Table1 and Table2) shifts up by one row.
$ws.Rows.Item(1).Delete()

# Remove the two stray blank rows that used to sit between the trailing
# note ("Senario") and the second table (rows 16-17 after the shift above).
[void]$ws.Range("A16:A17").EntireRow.Select()
$ws.Range("A16:A17").EntireRow.Delete()

# Set the sheet to print landscape (prepping the doc for a PDF export).
$ws.PageSetup.Orientation = 2
